$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update column widths (stored width = ColumnWidth + 5/6)
$ws.Columns.Item(2).ColumnWidth = 7.166666666666667
$ws.Columns.Item(3).ColumnWidth = 7.166666666666667
$ws.Columns.Item(6).ColumnWidth = 7.166666666666667
$ws.Columns.Item(7).ColumnWidth = 7.166666666666667
$ws.Columns.Item(9).ColumnWidth = 7.166666666666667
$ws.Columns.Item(10).ColumnWidth = 7.166666666666667
$ws.Columns.Item(11).ColumnWidth = 7.166666666666667
$ws.Columns.Item(12).ColumnWidth = 7.166666666666667
$ws.Columns.Item(13).ColumnWidth = 7.166666666666667
$ws.Columns.Item(15).ColumnWidth = 7.166666666666667
$ws.Columns.Item(16).ColumnWidth = 7.166666666666667
$ws.Columns.Item(17).ColumnWidth = 7.166666666666667
$ws.Columns.Item(20).ColumnWidth = 8.166666666666666
$ws.Columns.Item(23).ColumnWidth = 7.166666666666667
$ws.Columns.Item(24).ColumnWidth = 7.166666666666667
$ws.Columns.Item(26).ColumnWidth = 7.166666666666667
$ws.Columns.Item(27).ColumnWidth = 7.166666666666667
$ws.Columns.Item(28).ColumnWidth = 7.166666666666667
$ws.Columns.Item(29).ColumnWidth = 7.166666666666667
$ws.Columns.Item(30).ColumnWidth = 7.166666666666667
$ws.Columns.Item(34).ColumnWidth = 7.166666666666667

# Update data values in rows 2-5
# Row 2
$ws.Cells.Item(2, 1).Value = 45041.50694444445
$ws.Cells.Item(2, 2).Value = 4.639
$ws.Cells.Item(2, 3).Value = 3.949
$ws.Cells.Item(2, 4).Value = 0
$ws.Cells.Item(2, 5).Value = 5.682
$ws.Cells.Item(2, 6).Value = 5.56
$ws.Cells.Item(2, 7).Value = 1.404
$ws.Cells.Item(2, 8).Value = 7.057
$ws.Cells.Item(2, 9).Value = 3.269
$ws.Cells.Item(2, 10).Value = 2.974
$ws.Cells.Item(2, 11).Value = 2.59
$ws.Cells.Item(2, 12).Value = 3.719
$ws.Cells.Item(2, 13).Value = 4.974
$ws.Cells.Item(2, 14).Value = 2.211
$ws.Cells.Item(2, 15).Value = 2.2
$ws.Cells.Item(2, 16).Value = 3.582
$ws.Cells.Item(2, 17).Value = 1.706
$ws.Cells.Item(2, 18).Value = 0.605
$ws.Cells.Item(2, 19).Value = 0.024
$ws.Cells.Item(2, 20).Value = 37.258
$ws.Cells.Item(2, 21).Value = 6.598
$ws.Cells.Item(2, 22).Value = 4.255
$ws.Cells.Item(2, 23).Value = 5.158
$ws.Cells.Item(2, 24).Value = 1.722
$ws.Cells.Item(2, 25).Value = 0.377
$ws.Cells.Item(2, 26).Value = 3.036
$ws.Cells.Item(2, 27).Value = 1.789
$ws.Cells.Item(2, 28).Value = 1.745
$ws.Cells.Item(2, 29).Value = 5.904
$ws.Cells.Item(2, 30).Value = 3.659
$ws.Cells.Item(2, 31).Value = 3.978
$ws.Cells.Item(2, 32).Value = 5.202
$ws.Cells.Item(2, 33).Value = 0.995
$ws.Cells.Item(2, 34).Value = 2.83

# Row 3
$ws.Cells.Item(3, 1).Value = 45041.51388888889
$ws.Cells.Item(3, 2).Value = 13.336
$ws.Cells.Item(3, 3).Value = 10.154
$ws.Cells.Item(3, 4).Value = 0.273
$ws.Cells.Item(3, 5).Value = 27.215
$ws.Cells.Item(3, 6).Value = 22.999
$ws.Cells.Item(3, 7).Value = 9.476000000000001
$ws.Cells.Item(3, 8).Value = 34.049
$ws.Cells.Item(3, 9).Value = 15.211
$ws.Cells.Item(3, 10).Value = 7.641
$ws.Cells.Item(3, 11).Value = 10.427
$ws.Cells.Item(3, 12).Value = 11.395
$ws.Cells.Item(3, 13).Value = 12.478
$ws.Cells.Item(3, 14).Value = 3.868
$ws.Cells.Item(3, 15).Value = 9.768000000000001
$ws.Cells.Item(3, 16).Value = 14.375
$ws.Cells.Item(3, 17).Value = 7.92
$ws.Cells.Item(3, 18).Value = 0.554
$ws.Cells.Item(3, 19).Value = 0.286
$ws.Cells.Item(3, 20).Value = 146.748
$ws.Cells.Item(3, 21).Value = 27.605
$ws.Cells.Item(3, 22).Value = 9.861000000000001
$ws.Cells.Item(3, 23).Value = 19.024
$ws.Cells.Item(3, 24).Value = 9.547000000000001
$ws.Cells.Item(3, 25).Value = 1.328
$ws.Cells.Item(3, 26).Value = 17.462
$ws.Cells.Item(3, 27).Value = 8.081
$ws.Cells.Item(3, 28).Value = 7.29
$ws.Cells.Item(3, 29).Value = 9.797000000000001
$ws.Cells.Item(3, 30).Value = 11.983
$ws.Cells.Item(3, 31).Value = 1.327
$ws.Cells.Item(3, 32).Value = 30.593
$ws.Cells.Item(3, 33).Value = 5.084
$ws.Cells.Item(3, 34).Value = 11.39

# Row 4
$ws.Cells.Item(4, 1).Value = 45041.52083333334
$ws.Cells.Item(4, 2).Value = 20.894
$ws.Cells.Item(4, 3).Value = 15.773
$ws.Cells.Item(4, 4).Value = 0.5570000000000001
$ws.Cells.Item(4, 5).Value = 44.334
$ws.Cells.Item(4, 6).Value = 36.965
$ws.Cells.Item(4, 7).Value = 15.785
$ws.Cells.Item(4, 8).Value = 61.076
$ws.Cells.Item(4, 9).Value = 24.722
$ws.Cells.Item(4, 10).Value = 11.659
$ws.Cells.Item(4, 11).Value = 16.738
$ws.Cells.Item(4, 12).Value = 18.063
$ws.Cells.Item(4, 13).Value = 19.372
$ws.Cells.Item(4, 14).Value = 5.604
$ws.Cells.Item(4, 15).Value = 15.919
$ws.Cells.Item(4, 16).Value = 23.078
$ws.Cells.Item(4, 17).Value = 13.035
$ws.Cells.Item(4, 18).Value = 0.511
$ws.Cells.Item(4, 19).Value = 0.522
$ws.Cells.Item(4, 20).Value = 238.498
$ws.Cells.Item(4, 21).Value = 44.758
$ws.Cells.Item(4, 22).Value = 15.21
$ws.Cells.Item(4, 23).Value = 30.517
$ws.Cells.Item(4, 24).Value = 15.739
$ws.Cells.Item(4, 25).Value = 2.118
$ws.Cells.Item(4, 26).Value = 30.282
$ws.Cells.Item(4, 27).Value = 13.085
$ws.Cells.Item(4, 28).Value = 11.636
$ws.Cells.Item(4, 29).Value = 14.377
$ws.Cells.Item(4, 30).Value = 18.999
$ws.Cells.Item(4, 31).Value = 0.796
$ws.Cells.Item(4, 32).Value = 55.495
$ws.Cells.Item(4, 33).Value = 8.333
$ws.Cells.Item(4, 34).Value = 18.438

# Row 5
$ws.Cells.Item(5, 1).Value = 45041.52777777778
$ws.Cells.Item(5, 2).Value = 8.82
$ws.Cells.Item(5, 3).Value = 6.68
$ws.Cells.Item(5, 4).Value = 0.18
$ws.Cells.Item(5, 5).Value = 18.38
$ws.Cells.Item(5, 6).Value = 15.41
$ws.Cells.Item(5, 7).Value = 6.46
$ws.Cells.Item(5, 8).Value = 30.48
$ws.Cells.Item(5, 9).Value = 10.27
$ws.Cells.Item(5, 10).Value = 5.05
$ws.Cells.Item(5, 11).Value = 7
$ws.Cells.Item(5, 12).Value = 7.57
$ws.Cells.Item(5, 13).Value = 8.19
$ws.Cells.Item(5, 14).Value = 2.49
$ws.Cells.Item(5, 15).Value = 6.59
$ws.Cells.Item(5, 16).Value = 9.699999999999999
$ws.Cells.Item(5, 17).Value = 5.33
$ws.Cells.Item(5, 18).Value = 0.34
$ws.Cells.Item(5, 19).Value = 0.18
$ws.Cells.Item(5, 20).Value = 95.52
$ws.Cells.Item(5, 21).Value = 18.77
$ws.Cells.Item(5, 22).Value = 6.45
$ws.Cells.Item(5, 23).Value = 12.86
$ws.Cells.Item(5, 24).Value = 6.51
$ws.Cells.Item(5, 25).Value = 0.89
$ws.Cells.Item(5, 26).Value = 14.45
$ws.Cells.Item(5, 27).Value = 5.46
$ws.Cells.Item(5, 28).Value = 4.92
$ws.Cells.Item(5, 29).Value = 6.25
$ws.Cells.Item(5, 30).Value = 7.96
$ws.Cells.Item(5, 31).Value = 0.57
$ws.Cells.Item(5, 32).Value = 27.92
$ws.Cells.Item(5, 33).Value = 3.43
$ws.Cells.Item(5, 34).Value = 7.65

# Delete row 6 (data reduced from 6 rows to 5 rows)
$ws.Rows.Item(6).Delete()
